$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999477235368,
    0.99432659191773864,
    0.97259581828992336,
    0.96243797729675817,
    0.95267837777738729,
    0.92922202534380938,
    0.92261955021444875,
    0.91348256992667554,
    0.90172797308962471,
    0.89089337210757358,
    0.88929651352782468,
    0.88651545659602393,
    0.87522417742248149,
    0.87105587931224715,
    0.86846391298761894,
    0.86595708249950554,
    0.86224881171565937,
    0.86113987551425109,
    0.99127875826409428,
    0.96643503740420278,
    0.95972946962965366,
    0.95846494650669922,
    0.96556014514265454,
    0.95253837637782701,
    0.94608122874010014,
    0.9279195322337993,
    0.92702744212138688,
    0.92330298032901903,
    0.92128905272462103,
    0.92159630045876728,
    0.93076181270434666,
    0.93487631839238583,
    0.93435630503346534
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

